# [AA | 16/04/2018] : Changes for Profile tab and Framework restruture
#
# Adds three new front-of-book sheets (Visualize, DescriptiveStatistics_Data,
# QuantileStatistics_Data) ahead of the existing Categorical_Data / Profile_Data /
# Sample_Data / Fields_Data / Sheet1 sheets, populates them, and makes
# DescriptiveStatistics_Data the active tab.

$wb = $excel.ActiveWorkbook

$firstExisting = $wb.Worksheets.Item(1)

# Create in this order so the auto-assigned sheetId sequence comes out as
# Visualize=6, QuantileStatistics_Data=7, DescriptiveStatistics_Data=8.
$sVisualize = $wb.Worksheets.Add($firstExisting)
$sVisualize.Name = "Visualize"

$sQuantile = $wb.Worksheets.Add($firstExisting)
$sQuantile.Name = "QuantileStatistics_Data"

$sDescriptive = $wb.Worksheets.Add($firstExisting)
$sDescriptive.Name = "DescriptiveStatistics_Data"

# Re-fetch by name and reorder into the final tab order:
# Visualize, DescriptiveStatistics_Data, QuantileStatistics_Data, Categorical_Data, ...
$wsVisualize = $wb.Worksheets.Item("Visualize")
$wsDescriptive = $wb.Worksheets.Item("DescriptiveStatistics_Data")
$wsQuantile = $wb.Worksheets.Item("QuantileStatistics_Data")

[void]$wsVisualize.Move($wsDescriptive)
[void]$wsDescriptive.Move($wsQuantile)

# Moving a sheet invalidates previously-held references in this host, so
# re-resolve fresh handles by name before writing any cell data to them.
$wsVisualize = $wb.Worksheets.Item("Visualize")
$wsDescriptive = $wb.Worksheets.Item("DescriptiveStatistics_Data")
$wsQuantile = $wb.Worksheets.Item("QuantileStatistics_Data")

# --- Visualize: ID / Name / HireDate / Job / City -------------------------
$visRows = @(
    @(104, "AAVESH",   29923, "ANALYST",   "DALLAS"),
    @(100, "AMAN",     31920, "PRESIDENT", "DALLAS"),
    @(105, "SANTOSH",  29923, "CLERK",     "CHICAGO"),
    @(109, "MILLER",   29974, "CLERK",     "NEW YORK"),
    @(103, "SUYOG",    29746, "ANALYST",   "NEW YORK"),
    @(106, "SANGEETA", 29678, "MANAGER",   "DALLAS"),
    @(101, "GIRISH",   29637, "ANALYST",   "CHICAGO"),
    @(102, "AMOL",     29707, "MANAGER",   "CHICAGO"),
    @(107, "NISHA",    29907, "MANAGER",   "NEW YORK"),
    @(108, "MARTIN",   29857, "SALESMAN",  "CHICAGO"),
    @(110, "SCOTT",    31886, "ANALYST",   "DALLAS"),
    @(110, "SCOTT",    31886, "ANALYST",   "DALLAS")
)

$r = 1
foreach ($row in $visRows) {
    $wsVisualize.Cells.Item($r, 1).Value = $row[0]
    $wsVisualize.Cells.Item($r, 2).Value = $row[1]
    $wsVisualize.Cells.Item($r, 3).Value = $row[2]
    $wsVisualize.Cells.Item($r, 3).NumberFormat = "mm-dd-yy"
    $wsVisualize.Cells.Item($r, 4).Value = $row[3]
    $wsVisualize.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$wsVisualize.Columns.Item(3).ColumnWidth = 10.71
$wsVisualize.Columns.Item(4).ColumnWidth = 10.57
$wsVisualize.Columns.Item(5).ColumnWidth = 10.43

# --- DescriptiveStatistics_Data --------------------------------------------
$descRows = @(
    @("Standard deviation", 226.58),
    @("Coef of variation",  0.54),
    @("Kurtosis",          -1.05),
    @("Mean",               421.89),
    @("MAD",                185.01),
    @("Skewness",           0.13),
    @("Sum",                3797),
    @("Variance",           51338.36),
    @("Memory size",        3560)
)

$r = 1
foreach ($row in $descRows) {
    $wsDescriptive.Cells.Item($r, 1).Value = $row[0]
    $wsDescriptive.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
$wsDescriptive.Columns.Item(1).ColumnWidth = 18
[void]$wsDescriptive.Range("A1:B9").Select()

# --- QuantileStatistics_Data ------------------------------------------------
$quantRows = @(
    @("Minimum",              100),
    @("5-th Percentile",      100),
    @("Q1",                   211),
    @("Median",                444),
    @("Q3",                   610.5),
    @("95-th Percentile",     777),
    @("Maximum",               777),
    @("Range",                677),
    @("Interquartile range",  399.5)
)

$r = 1
foreach ($row in $quantRows) {
    $wsQuantile.Cells.Item($r, 1).Value = $row[0]
    $wsQuantile.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
$wsQuantile.Columns.Item(1).ColumnWidth = 18
$wsQuantile.Columns.Item(2).ColumnWidth = 6
[void]$wsQuantile.Range("A1:B9").Select()

# Make DescriptiveStatistics_Data the active tab (activeTab/firstSheet = 1)
$wb.Worksheets.Item("DescriptiveStatistics_Data").Activate()

Write-Host "done"

